# Update TPM-derived values on Sheet1 (Rspo1-Lrp6) to reflect new TPM recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 13.17295566666667
$ws.Range("N2").Value = 39.518867
$ws.Range("O2").Value = 0.133784132206724
$ws.Range("P2").Value = 0.133784132206724
$ws.Range("Q2").Value = 2.077001874889445
$ws.Range("R2").Value = 18.693016874005
$ws.Range("S2").Value = 0.133784132206724
$ws.Range("T2").Value = 0.133784132206724

# Row 3 (M3/N3 unchanged)
$ws.Range("O3").Value = 0.4382627974978752
$ws.Range("P3").Value = 0.4382627974978752
$ws.Range("Q3").Value = 6.804040487333889
$ws.Range("R3").Value = 61.236364386005
$ws.Range("S3").Value = 0.4382627974978752
$ws.Range("T3").Value = 0.4382627974978752

# Row 4
$ws.Range("M4").Value = 21.06166566666667
$ws.Range("N4").Value = 63.184997
$ws.Range("O4").Value = 0.2139016281041017
$ws.Range("P4").Value = 0.2139016281041017
$ws.Range("Q4").Value = 3.320827928439445
$ws.Range("R4").Value = 29.887451355955
$ws.Range("S4").Value = 0.2139016281041017
$ws.Range("T4").Value = 0.2139016281041017

# Row 5
$ws.Range("M5").Value = 21.076417
$ws.Range("N5").Value = 63.229251
$ws.Range("O5").Value = 0.214051442191299
$ws.Range("P5").Value = 0.214051442191299
$ws.Range("Q5").Value = 3.323153795751667
$ws.Range("R5").Value = 29.908384161765
$ws.Range("S5").Value = 0.214051442191299
$ws.Range("T5").Value = 0.214051442191299
